$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rundata")

# Update "CrepeErase" postpone days value (row 5, column D) from 20 to 31
$ws.Range("D5").Value = 31

# Remove the "SeaCalmSkin" row (row 7) entirely; rows below shift up
$ws.Rows.Item(7).Delete()

# Update the selection to match the resulting layout
$ws.Activate()
$ws.Range("E7").Select()
